# Fixed bug where 0 is first number at checkout
#
# Barcode 045496870775 ("Doritos") was scanning incorrectly at checkout
# whenever the leading 0 of the barcode was dropped. Renamed the inventory
# item to "Doritos2" to distinguish it during testing, logged the two test
# checkouts that exposed the issue, and corrected the remaining stock count.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # inventory
$ws2 = $wb.Worksheets.Item(2)   # transactions

# --- Rename the sheets to the new snapshot date -----------------------------
$ws1.Name = "inventory 15-02-2024"
$ws2.Name = "transactions 15-02-2024"

# --- Inventory: rename the item and correct its remaining quantity ---------
# (50 in stock - 2 - 2 consumed by the two test checkouts below = 46)
$ws1.Range("C3").Value = "Doritos2"
$ws1.Range("E3").Value = 46

# --- Transactions: log the two checkouts that exposed the bug --------------
# Seed rows 10/11 from the last existing row so number formats/styles (e.g.
# the barcode staying text, the bold/border style on column A) carry over
# exactly as they do for the existing rows, then overwrite the values.
$ws2.Range("A9:H9").Copy($ws2.Range("A10:H10"))
$ws2.Range("A9:H9").Copy($ws2.Range("A11:H11"))

$ws2.Range("A10").Value = 8
$ws1.Range("B3").Copy($ws2.Range("B10"))     # barcode 045496870775, kept as text
$ws2.Range("C10").Value = "Doritos2"
$ws2.Range("D10").Value = 2
$ws2.Range("E10").Value = 2
$ws2.Range("F10").Value = "Snacks"
$ws2.Range("G10").Value = "Costco"
$ws2.Range("H10").Value = "30-01-2024 15:12:14"

$ws2.Range("A11").Value = 9
$ws1.Range("B3").Copy($ws2.Range("B11"))     # barcode 045496870775, kept as text
$ws2.Range("C11").Value = "Doritos2"
$ws2.Range("D11").Value = 2
$ws2.Range("E11").Value = 2
$ws2.Range("F11").Value = "Snacks"
$ws2.Range("G11").Value = "Costco"
$ws2.Range("H11").Value = "30-01-2024 15:13:07"
